$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.412.59'
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").Value = '1.669.80'
$ws.Range("E3").Value = '  -4.03%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.68%  '
$ws.Range("D5").Value = '239.74'
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("D7").Value = '0.4829'
$ws.Range("E7").Value = '  -7.34%  '
$ws.Range("D8").Value = '0.2637'
$ws.Range("E8").Value = '  -4.41%  '
$ws.Range("D9").Value = '0.06015'
$ws.Range("E9").Value = '  -2.39%  '
$ws.Range("D10").Value = '0.07157'
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("D11").Value = '1.693.82'
$ws.Range("E11").Value = '  -2.67%  '
$ws.Range("D12").Value = '0.6277'
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("D13").Value = '14.55'
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("D14").Value = '4.655'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = '73.55'
$ws.Range("E15").Value = '  -5.32%  '
$ws.Range("B16").Value = 'Dai'
$ws.Range("C16").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = '25.420.35'
$ws.Range("E18").Value = '  -1.93%  '
$ws.Range("D19").Value = '11.58'
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").Value = '0.000006644'
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("D21").Value = '1.917.21'
$ws.Range("E21").Value = '  -2.49%  '
$ws.Range("D22").Value = '4.468'
$ws.Range("E22").Value = '  +4.11%  '
$ws.Range("D23").Value = '8.630'
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").Value = '5.310'
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("D25").Value = '134.11'
$ws.Range("E25").Value = '  -3.59%  '
$ws.Range("E26").Value = '  -2.30%  '
$ws.Range("D27").Value = '1.380'
$ws.Range("E27").Value = '  -9.55%  '
$ws.Range("D28").Value = '1.724'
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("D29").Value = '102.79'
$ws.Range("E29").Value = '  -3.17%  '
$ws.Range("D30").Value = '3.885'
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("D31").Value = '0.07960'
$ws.Range("E31").Value = '  -4.19%  '
$ws.Range("D32").Value = '3.553'
$ws.Range("E32").Value = '  -3.97%  '
$ws.Range("D33").Value = '0.04638'
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("D34").Value = '2.647'
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("D35").Value = '0.9538'
$ws.Range("E35").Value = '  -3.79%  '
$ws.Range("D36").Value = '0.5877'
$ws.Range("E36").Value = '  -5.29%  '
$ws.Range("D37").Value = '2.647'
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").Value = '0.01559'
$ws.Range("E38").Value = '  -3.12%  '
$ws.Range("D39").Value = '1.007'
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("D40").Value = '0.8404'
$ws.Range("E40").Value = '  +12.92%  '
$ws.Range("D41").Value = '1.879'
$ws.Range("D42").Value = '99.17'
$ws.Range("E42").Value = '  +1.45%  '
$ws.Range("D43").Value = '0.3761'
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("D44").Value = '4.910'
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("D45").Value = '0.1150'
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("D46").Value = '6.108'
$ws.Range("E46").Value = '  -2.47%  '
$ws.Range("D47").Value = '0.05192'
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("D48").Value = '53.97'
$ws.Range("E48").Value = '  -1.65%  '
$ws.Range("D49").Value = '29.80'
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("B50").Value = 'TrueUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D50").Value = '1.007'
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.373'
$ws.Range("E51").Value = '  -3.30%  '
